$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9568136666666667
$ws.Range("H2").Value = 2.870441
$ws.Range("I2").Value = 0.09967139189263423
$ws.Range("J2").Value = 0.09967139189263421
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 87.55926664062588
$ws.Range("R2").Value = 788.0333997656329
$ws.Range("S2").Value = 0.09653692417098364
$ws.Range("T2").Value = 0.09653692417098364
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9568136666666667
$ws.Range("H3").Value = 2.870441
$ws.Range("I3").Value = 0.09967139189263423
$ws.Range("J3").Value = 0.09967139189263421
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 0.128894282664
$ws.Range("R3").Value = 1.160048543976
$ws.Range("S3").Value = 0.0001421101165988358
$ws.Range("T3").Value = 0.0001421101165988358
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9568136666666667
$ws.Range("H4").Value = 2.870441
$ws.Range("I4").Value = 0.09967139189263423
$ws.Range("J4").Value = 0.09967139189263421
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 2.714076915903889
$ws.Range("R4").Value = 24.426692243135
$ws.Range("S4").Value = 0.00299235760505175
$ws.Range("T4").Value = 0.00299235760505175
$ws.Range("I5").Value = 0.789130862182032
$ws.Range("J5").Value = 0.789130862182032
$ws.Range("M5").Value = 91.51130433333333
$ws.Range("N5").Value = 274.533913
$ws.Range("O5").Value = 0.9685519820468944
$ws.Range("P5").Value = 0.9685519820468945
$ws.Range("Q5").Value = 693.2352229070231
$ws.Range("R5").Value = 6239.117006163208
$ws.Range("S5").Value = 0.7643142606607818
$ws.Range("T5").Value = 0.7643142606607819
$ws.Range("I6").Value = 0.789130862182032
$ws.Range("J6").Value = 0.789130862182032
$ws.Range("O6").Value = 0.001425786415744213
$ws.Range("P6").Value = 0.001425786415744214
$ws.Range("S6").Value = 0.00112513206354366
$ws.Range("T6").Value = 0.00112513206354366
$ws.Range("I7").Value = 0.789130862182032
$ws.Range("J7").Value = 0.789130862182032
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.03002223153736139
$ws.Range("P7").Value = 0.03002223153736139
$ws.Range("S7").Value = 0.02369146945770659
$ws.Range("T7").Value = 0.02369146945770659
$ws.Range("I8").Value = 0.1111977459253338
$ws.Range("J8").Value = 0.1111977459253338
$ws.Range("M8").Value = 91.51130433333333
$ws.Range("N8").Value = 274.533913
$ws.Range("O8").Value = 0.9685519820468944
$ws.Range("P8").Value = 0.9685519820468945
$ws.Range("Q8").Value = 97.684931457573
$ws.Range("R8").Value = 879.164383118157
$ws.Range("S8").Value = 0.107700797215129
$ws.Range("T8").Value = 0.107700797215129
$ws.Range("I9").Value = 0.1111977459253338
$ws.Range("J9").Value = 0.1111977459253338
$ws.Range("O9").Value = 0.001425786415744213
$ws.Range("P9").Value = 0.001425786415744214
$ws.Range("S9").Value = 0.0001585442356017174
$ws.Range("T9").Value = 0.0001585442356017174
$ws.Range("I10").Value = 0.1111977459253338
$ws.Range("J10").Value = 0.1111977459253338
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509734999999999
$ws.Range("O10").Value = 0.03002223153736139
$ws.Range("P10").Value = 0.03002223153736139
$ws.Range("Q10").Value = 3.027942417435
$ws.Range("R10").Value = 27.251481756915
$ws.Range("S10").Value = 0.003338404474603055
$ws.Range("T10").Value = 0.003338404474603055
